$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) Shift columns N(14)..AC(29) right by one column -> O(15)..AD(30)
#    Process from the right-most column down to the left-most so we never
#    clobber a source cell before it has been read.
for ($c = 29; $c -ge 14; $c--) {
    for ($r = 1; $r -le 127; $r++) {
        $v = $ws.Cells.Item($r, $c).Value()
        $ws.Cells.Item($r, $c + 1).Value = $v
    }
}


# 2) Write the brand-new "localdb" column (N) header + its six commands
$ws.Cells.Item(1, 14).Value = "localdb"
$ws.Cells.Item(2, 14).Value = "cloneTable(var,source,target)"
$ws.Cells.Item(3, 14).Value = "dropTables(var,tables)"
$ws.Cells.Item(4, 14).Value = "exportCSV(sql,output)"
$ws.Cells.Item(5, 14).Value = "importRecords(var,sourceDb,sql,table)"
$ws.Cells.Item(6, 14).Value = "purge(var)"
$ws.Cells.Item(7, 14).Value = "runSQLs(var,sqls)"


# 3) Insert "localdb" into the alphabetical "target" list in column A,
#    between "json" (row 13) and "macro" (row 14), pushing the remainder
#    down by one row.
for ($r = 29; $r -ge 14; $r--) {
    $v = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r + 1, 1).Value = $v
}
$ws.Cells.Item(14, 1).Value = "localdb"

# 4) Re-point every defined name whose target column moved right by one
#    (everything from the old "mail" name at column O onward), widen
#    "target" by the extra row, and add the brand-new "localdb" name.
$wb.Names.Item("mail").RefersTo      = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo    = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo       = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo     = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo     = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo       = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo     = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo       = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo      = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo    = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo       = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo  = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo        = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo  = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo       = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("macro").RefersTo     = "='#system'!`$O`$2:`$O`$4"
$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")

Write-Output "shift done"
